$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the "Last Updated" timestamp by one minute ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "05 Nov 2025, 01:08 PM"

# --- Stock List sheet: a new entry (CAPTRU-RE1) is published at the top of
#     the list. Every existing row shifts down by one, and the row that used
#     to be last (TRAVELFOOD, row 76) falls off the bottom since the sheet
#     keeps a fixed number of rows. Columns A (icon), F and G (both always
#     "N/A" in this sheet) are untouched by the shift.
$ws = $wb.Worksheets.Item("Stock List")

for ($r = 76; $r -ge 3; $r--) {
    $src = $r - 1
    $ws.Range("B$r").Value = $ws.Range("B$src").Value2
    $ws.Range("C$r").Value = $ws.Range("C$src").Value2
    $ws.Range("D$r").Value = $ws.Range("D$src").Value2
    $ws.Range("E$r").Value = $ws.Range("E$src").Value2
    $ws.Range("H$r").Value = $ws.Range("H$src").Value2
}

$ws.Range("B2").Value = "CAPTRU-RE1"
$ws.Range("C2").Value = "CAPTRU-RE1"
$ws.Range("D2").Value = 5.67
$ws.Range("E2").Value = -11.9565
$ws.Range("H2").Value = 0
